$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.328.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.488.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0810'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.22%  '
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.878.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.502.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.845'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.230.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("E22").Value = '  +13.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '245.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.20%  '
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.136'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.00%  '
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.66'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.88%  '
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '118.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.37%  '
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.984.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  -6.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("E50").Value = '  -6.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.77%  '
